# Add new columns I (I0) and J (IF) to Sheet1, mirroring style/format of
# existing header (H1) and data columns (H2:H58).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Match the styling used by the other header cells (e.g. H1: bold font,
# thin border, centered alignment) by copying its format onto I1:J1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data values for column I (rows 2-58) ---
$iValues = @(8, 7, 5, 8, 8, 6, 7, 6, 7, 8, 7, 6, 7, 8, 7, 7, 9, 5, 6, 8, 8, 7, 7, 4, 6, 6, 9, 8, 6, 6, 10, 7, 9, 8, 5, 7, 6, 6, 8, 14, 8, 9, 10, 8, 7, 10, 5, 9, 9, 8, 8, 5, 8, 6, 5, 7, 6)

# --- Data values for column J (rows 2-58) ---
$jValues = @(8, 8, 6, 8, 8, 6, 7, 7, 7, 8, 7, 6, 7, 8, 7, 7, 9, 7, 7, 8, 8, 8, 7, 5, 6, 7, 9, 8, 6, 7, 10, 8, 9, 9, 6, 8, 6, 8, 8, 14, 8, 9, 10, 8, 8, 10, 5, 9, 9, 8, 8, 5, 8, 6, 5, 7, 6)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
